$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.224.50"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "1.843.74"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.38"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6276"
$ws.Range("E6").Value = "  -1.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07524"
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2949"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("E10").Value = "  +1.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07725"
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("D12").Value = "1.875.34"
$ws.Range("E12").Value = "  +1.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.038"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6787"
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.21"
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009241"
$ws.Range("E16").Value = "  -3.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.989"
$ws.Range("E17").Value = "  -2.20%  "
$ws.Range("D18").Value = "29.241.60"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").Value = "2.111.67"
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "232.10"
$ws.Range("E20").Value = "  +2.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.74"
$ws.Range("E21").Value = "  +0.90%  "
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.53"
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1397"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.563"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.96"
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.200"
$ws.Range("E30").Value = "  +1.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.167"
$ws.Range("E31").Value = "  +1.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05580"
$ws.Range("E32").Value = "  +3.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.209"
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7529"
$ws.Range("E34").Value = "  +0.47%  "
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("E36").Value = "  +0.56%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").Value = "1.239.02"
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.771"
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01797"
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.638"
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9019"
$ws.Range("E42").Value = "  -0.61%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.52"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").Value = "1.995.71"
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.72"
$ws.Range("E46").Value = "  +2.19%  "
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("E48").Value = "  -4.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4104"
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.129"
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05848"
$ws.Range("E51").Value = "  +1.14%  "
